# Trade #192 (0-indexed trade #220) closed at 2026-02-17 22:07:51 - unknown UNKNOWN +0.000%
# Applies:
#   - Summary sheet roll-up numbers
#   - Strategy Status roll-up numbers for MarketMaking
#   - Closes the open MarketMaking trade recorded at All Trades row 221 / MarketMaking row 188
#   - Appends two freshly-opened trades (volatility_scorer #253, MarketMaking #254)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1399.97   # Current Capital
$summary.Range("B4").Value = -0.25     # Total P&L $
$summary.Range("B5").Value = -0.02     # Total P&L %
$summary.Range("B6").Value = 220       # Total Trades
$summary.Range("B7").Value = 85        # Winning Trades
$summary.Range("B9").Value = 38.64     # Win Rate %

# ---------------------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 99.97   # Capital
$status.Range("D5").Value = 187     # Trades
$status.Range("E5").Value = -0.36   # P&L $
$status.Range("F5").Value = -0.03   # P&L %
$status.Range("G5").Value = 37.97   # Win Rate %

# ---------------------------------------------------------------------------
# 3. All Trades sheet - close trade #220 (row 221)
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Cells.Item(221, 7).Value = 0.31          # G Exit Price
$allTrades.Cells.Item(221, 8).Value = "CLOSED"      # H Status
$allTrades.Cells.Item(221, 9).Value = 40.9091       # I P&L %
$allTrades.Cells.Item(221, 10).Value = 0.09         # J P&L $
$allTrades.Cells.Item(221, 11).Value = 99.97        # K Capital After
$allTrades.Cells.Item(221, 12).Value = "early_exit" # L Exit Reason
$allTrades.Cells.Item(221, 13).Value = 0.18         # M Duration (min)

# ---------------------------------------------------------------------------
# 4. MarketMaking sheet - close trade #220 (row 188) - same trade, different
#    column layout (Exit Reason/Duration are columns P/Q here)
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Cells.Item(188, 7).Value = 0.31          # G Exit Price
$marketMaking.Cells.Item(188, 8).Value = "CLOSED"      # H Status
$marketMaking.Cells.Item(188, 9).Value = 40.9091       # I P&L %
$marketMaking.Cells.Item(188, 10).Value = 0.09         # J P&L $
$marketMaking.Cells.Item(188, 11).Value = 99.97        # K Capital After
$marketMaking.Cells.Item(188, 16).Value = "early_exit" # P Exit Reason
$marketMaking.Cells.Item(188, 17).Value = 0.18         # Q Duration (min)

# ---------------------------------------------------------------------------
# 5. All Trades sheet - append newly opened trades #253 and #254
# ---------------------------------------------------------------------------
function Set-TextCell($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# Row 254 -> trade #253 (volatility_scorer, still OPEN)
Set-TextCell $allTrades 254 2 "2026-02-17"
Set-TextCell $allTrades 254 3 "22:07:44"
$allTrades.Cells.Item(254, 1).Value = 253
$allTrades.Cells.Item(254, 4).Value = "volatility_scorer"
$allTrades.Cells.Item(254, 5).Value = "NEUTRAL"
$allTrades.Cells.Item(254, 6).Value = 0.22
$allTrades.Cells.Item(254, 8).Value = "OPEN"
$allTrades.Cells.Item(254, 9).Value = 0
$allTrades.Cells.Item(254, 10).Value = 0
$allTrades.Cells.Item(254, 11).Value = 100
$allTrades.Cells.Item(254, 13).Value = 0
$allTrades.Cells.Item(254, 14).Value = 0
$allTrades.Cells.Item(254, 15).Value = 0
$allTrades.Cells.Item(254, 16).Value = 0.85
$allTrades.Cells.Item(254, 17).Value = "Low vol market (score: inf) - ideal for market making"

# Row 255 -> trade #254 (MarketMaking, still OPEN)
Set-TextCell $allTrades 255 2 "2026-02-17"
Set-TextCell $allTrades 255 3 "22:07:44"
$allTrades.Cells.Item(255, 1).Value = 254
$allTrades.Cells.Item(255, 4).Value = "MarketMaking"
$allTrades.Cells.Item(255, 5).Value = "UP"
$allTrades.Cells.Item(255, 6).Value = 0.78
$allTrades.Cells.Item(255, 8).Value = "OPEN"
$allTrades.Cells.Item(255, 9).Value = 0
$allTrades.Cells.Item(255, 10).Value = 0
$allTrades.Cells.Item(255, 11).Value = 99.87837680355362
$allTrades.Cells.Item(255, 13).Value = 0
$allTrades.Cells.Item(255, 14).Value = 0
$allTrades.Cells.Item(255, 15).Value = 0
$allTrades.Cells.Item(255, 16).Value = 0.6
$allTrades.Cells.Item(255, 17).Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# 6. volatility_scorer sheet - append trade #253 (row 9)
# ---------------------------------------------------------------------------
$volScorer = $wb.Worksheets.Item("volatility_scorer")
Set-TextCell $volScorer 9 2 "2026-02-17"
Set-TextCell $volScorer 9 3 "22:07:44"
$volScorer.Cells.Item(9, 1).Value = 253
$volScorer.Cells.Item(9, 4).Value = "volatility_scorer"
$volScorer.Cells.Item(9, 5).Value = "NEUTRAL"
$volScorer.Cells.Item(9, 6).Value = 0.22
$volScorer.Cells.Item(9, 8).Value = "OPEN"
$volScorer.Cells.Item(9, 9).Value = 0
$volScorer.Cells.Item(9, 10).Value = 0
$volScorer.Cells.Item(9, 11).Value = 100
$volScorer.Cells.Item(9, 12).Value = 0
$volScorer.Cells.Item(9, 13).Value = 0
$volScorer.Cells.Item(9, 14).Value = 0.85
$volScorer.Cells.Item(9, 15).Value = "Low vol market (score: inf) - ideal for market making"
$volScorer.Cells.Item(9, 17).Value = 0

# ---------------------------------------------------------------------------
# 7. MarketMaking sheet - append trade #254 (row 214)
# ---------------------------------------------------------------------------
Set-TextCell $marketMaking 214 2 "2026-02-17"
Set-TextCell $marketMaking 214 3 "22:07:44"
$marketMaking.Cells.Item(214, 1).Value = 254
$marketMaking.Cells.Item(214, 4).Value = "MarketMaking"
$marketMaking.Cells.Item(214, 5).Value = "UP"
$marketMaking.Cells.Item(214, 6).Value = 0.78
$marketMaking.Cells.Item(214, 8).Value = "OPEN"
$marketMaking.Cells.Item(214, 9).Value = 0
$marketMaking.Cells.Item(214, 10).Value = 0
$marketMaking.Cells.Item(214, 11).Value = 99.87837680355362
$marketMaking.Cells.Item(214, 12).Value = 0
$marketMaking.Cells.Item(214, 13).Value = 0
$marketMaking.Cells.Item(214, 14).Value = 0.6
$marketMaking.Cells.Item(214, 15).Value = "Normal spread capture: 19600 bps"
$marketMaking.Cells.Item(214, 17).Value = 0

Write-Output "edit.ps1 applied"
